# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the underlying recalculation replaced the values in
# column G ("K") for every data row (rows 2-43) on the active sheet.
# We reproduce the resulting literal values directly, as shown by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 4
    10 = 3
    11 = 1
    12 = 4
    13 = 1
    14 = 2
    15 = 6
    16 = 4
    17 = 3
    18 = 3
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 3
    30 = 5
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 4
    37 = 2
    38 = 2
    39 = 0
    40 = 2
    41 = 1
    42 = 2
    43 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
